$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.1806205
$ws.Range("H2").Value = 0.361241
$ws.Range("I2").Value = 0.7284509268949775
$ws.Range("J2").Value = 0.7284509268949775
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.006466
$ws.Range("N2").Value = 0.012932
$ws.Range("O2").Value = 0.008493427970384656
$ws.Range("P2").Value = 0.008493427970384656
$ws.Range("Q2").Value = 0.001167892153
$ws.Range("R2").Value = 0.004671568612
$ws.Range("S2").Value = 0.00618704547754243
$ws.Range("T2").Value = 0.00618704547754243

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.1806205
$ws.Range("H3").Value = 0.361241
$ws.Range("I3").Value = 0.7284509268949775
$ws.Range("J3").Value = 0.7284509268949775
$ws.Range("M3").Value = 0.7548284999999999
$ws.Range("N3").Value = 1.509657
$ws.Range("O3").Value = 0.9915065720296153
$ws.Range("P3").Value = 0.9915065720296153
$ws.Range("Q3").Value = 0.13633750108425
$ws.Range("R3").Value = 0.5453500043369999
$ws.Range("S3").Value = 0.722263881417435
$ws.Range("T3").Value = 0.722263881417435

# Row 4 (new): MuSCs -> ECs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.067331
$ws.Range("H4").Value = 0.134662
$ws.Range("I4").Value = 0.2715490731050226
$ws.Range("J4").Value = 0.2715490731050226
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.006466
$ws.Range("N4").Value = 0.012932
$ws.Range("O4").Value = 0.008493427970384656
$ws.Range("P4").Value = 0.008493427970384656
$ws.Range("Q4").Value = 0.0004353622460000001
$ws.Range("R4").Value = 0.001741448984
$ws.Range("S4").Value = 0.002306382492842226
$ws.Range("T4").Value = 0.002306382492842226

# Row 5 (new): MuSCs -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.067331
$ws.Range("H5").Value = 0.134662
$ws.Range("I5").Value = 0.2715490731050226
$ws.Range("J5").Value = 0.2715490731050226
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7548284999999999
$ws.Range("N5").Value = 1.509657
$ws.Range("O5").Value = 0.9915065720296153
$ws.Range("P5").Value = 0.9915065720296153
$ws.Range("Q5").Value = 0.0508233577335
$ws.Range("R5").Value = 0.203293430934
$ws.Range("S5").Value = 0.2692426906121804
$ws.Range("T5").Value = 0.2692426906121804
